$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 172, shifting the
# existing rows 172:193 down to 174:195 (new weekly data is prepended
# to this block of "Vega Monumental Concepción" / "Zanahoria" records).
$ws.Rows("172:173").Insert()

# Fill in the two newly inserted rows with the new week's data. The
# surrounding (now shifted) rows keep the same constant columns
# (A, B, C, E, F, G, H, N, Q, R), so replicate those here too.

# Row 172 - "Primera" quality
$ws.Cells.Item(172, 1).Value = 11
$ws.Cells.Item(172, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(172, 3).Value = "Bíobío"
$ws.Cells.Item(172, 4).Value = 44617
$ws.Cells.Item(172, 5).Value = 8
$ws.Cells.Item(172, 6).Value = 100114013
$ws.Cells.Item(172, 7).Value = "Zanahoria"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 500
$ws.Cells.Item(172, 11).Value = 8000
$ws.Cells.Item(172, 12).Value = 8000
$ws.Cells.Item(172, 13).Value = 8000
$ws.Cells.Item(172, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(172, 15).Value = "Región de Ñuble"
$ws.Cells.Item(172, 16).Value = 400
$ws.Cells.Item(172, 17).Value = 20
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# Row 173 - "Segunda" quality
$ws.Cells.Item(173, 1).Value = 11
$ws.Cells.Item(173, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(173, 3).Value = "Bíobío"
$ws.Cells.Item(173, 4).Value = 44617
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = 100114013
$ws.Cells.Item(173, 7).Value = "Zanahoria"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Segunda"
$ws.Cells.Item(173, 10).Value = 500
$ws.Cells.Item(173, 11).Value = 7000
$ws.Cells.Item(173, 12).Value = 7000
$ws.Cells.Item(173, 13).Value = 7000
$ws.Cells.Item(173, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(173, 15).Value = "Región de Ñuble"
$ws.Cells.Item(173, 16).Value = 350
$ws.Cells.Item(173, 17).Value = 20
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Match the date formatting style already used in column D.
$ws.Range("D172:D173").NumberFormat = $ws.Range("D174").NumberFormat
